$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the CasesTab query text (B2): drop the trailing Cohort column line ---
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Female'] and demo.neutered_indicator IN ['Yes']
    
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value = $newCasesQuery

# --- Row heights (auto-adjusted by Excel after the text edit / app upgrade) ---
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 275.5
$ws.Rows.Item(4).RowHeight = 275.5

# --- Column widths (slightly narrower after re-save with newer Excel build) ---
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 38.833333

# --- View state: window scrolled up to row 2, selection moved to B2 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("B2").Select()

Write-Host "Edit applied."
